$d = $word.ActiveDocument

# Locate the paragraph containing the erratum ("canal 7" -> "canal 12").
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*del canal 7 del reservorio multicanal*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph containing 'del canal 7 del reservorio multicanal'."
}

$full = $target.Range.Text
$old  = " (x8) del canal 7 del reservorio multicanal a cada uno de los pocillos del deepwell del slot "
$idx  = $full.IndexOf($old)
if ($idx -lt 0) {
    throw "Could not locate the exact erratum text within the target paragraph."
}

$s = $target.Range.Start + $idx

$piece1 = " (x8) del canal "
$piece2 = "12"
$piece3 = " del reservorio multicanal a cada uno de los pocillos del deepwell del slot "
$newText = $piece1 + $piece2 + $piece3

# Replace the old run's text (" ... canal 7 ... slot ") with the corrected text.
$rng = $d.Range($s, $s + $old.Length)
$rng.Text = $newText

# Work out the character boundaries of the three new pieces.
$p1Start = $s
$p1End   = $p1Start + $piece1.Length
$p2Start = $p1End
$p2End   = $p2Start + $piece2.Length
$p3Start = $p2End
$p3End   = $p3Start + $piece3.Length

# The "4" (slot number) run that follows immediately must stay its own run too.
$fourStart = $p3End
$fourEnd   = $fourStart + 1

# Force run boundaries to match the original diff's run split: toggling a
# character formatting property on/off causes the engine to split the
# surrounding text into distinct runs at the touched boundaries, without
# altering the visible formatting (identical before/after).
$b1 = $d.Range($p1Start, $p1End)
$b1.Bold = 1
$b1.Bold = 0

$b2 = $d.Range($p2Start, $p2End)
$b2.Bold = 1
$b2.Bold = 0

$b3 = $d.Range($p3Start, $fourEnd)
$b3.Bold = 1
$b3.Bold = 0

$b4 = $d.Range($fourStart, $fourEnd)
$b4.Bold = 1
$b4.Bold = 0
